# Applies the "Added Introduction and help in Config file" edit to Config.xlsx
# Target: xl/worksheets/sheet1.xml (Introduction sheet) gains a new block of
# rows (Introduction text + per-sheet descriptions) inserted before the old
# "Legend" block, plus some cosmetic selection / column width changes.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Introduction")

# ---------------------------------------------------------------------------
# 1. Make room: insert 17 new blank rows right after row 2 (i.e. before the
#    old row 3), shifting the old rows 3-8 down to 20-25.
# ---------------------------------------------------------------------------
$ws.Range("A3:A19").EntireRow.Insert()

# ---------------------------------------------------------------------------
# 2. Set all the new cell text, in the same order the strings were originally
#    authored (section headings first, then descriptions) so that the shared
#    string table is built up in the expected order.
# ---------------------------------------------------------------------------
$ws.Range("A3").Value2 = "This is the configuration file used to describe various changeable parameters of the process. You should use this file to store settings that are environment related (like paths to programs or resources), user related (email account names, credential names), or plain data (URL of website or name of SAP report to execute). Below, the purpose of each sheet is explained in more detail."
$ws.Range("A4").Value2 = "Settings"
$ws.Range("A6").Value2 = "Credentials"
$ws.Range("A8").Value2 = "Workblocks"
$ws.Range("A10").Value2 = "Tasks"
$ws.Range("A15").Value2 = "Constants"
$ws.Range("A17").Value2 = "Assets"
$ws.Range("A7").Value2 = "The credentials sheet is the place to store your credential names.There is also one special credential, that needs to be defined only once, and which is comprised of the URL, TenancyName and CredentialName required to authenticate to the Orchestrator server using REST API. This is only used when working with QueueItems."
$ws.Range("A9").Value2 = "The workblock names are of the states in the framework. Define the names of workblocks you create here."
$ws.Range("A12").Value2 = "FirstRunTask: This task is invoked in the Framework Layer and executes only once(Even if Transaction number 1 were to fail and be retried, it would not be executed again), at program startup. It should not interact with data in memory, since it executes before we enter the main process data layer, but it can be used as a queue dispatcher."
$ws.Range("A14").Value2 = "Task1: This task is not invoked anywhere, and should be used by the developer. "
$ws.Range("A11").Value2 = "The sheet contains the list of tasks. Each task is another Business Process Layer context that is executed at some point during the main process execution. For system tasks, the execution is preselected and configurable from the settings. For user added tasks, it is chosen by the user."
$ws.Range("A13").Value2 = "GetDataTask: This task is invoked in the Data Layer of the main task. The reason is that we might desire it to deliver some TransactionData to us is a safe manner. Thus, it might navigate a website, download a file, process it, and deliver us an output datatable TransactionData. This would be made available in the Data Layer of the main task and would be ready for usage according to the business rules of the process."
$ws.Range("A5").Value2 = "This sheet is the place to store plain data, as well as most user data with the important exception of credential names."
$ws.Range("A16").Value2 = "Typically there is not much for you to add here, although you want to check/change the settings of the Retry mechanism implemented in at the framework layer, during transaction processing, exception recovery, and continuous failiure. Also stores constants used throughout the program, like preconfiguered delays, timeouts. "
$ws.Range("A18").Value2 = "This sheet is used to fetch assets from Orchestrator. The column name is the key, while the column asset hoolds the asset name in Orchestrator. If there is another local key with the same name, it will be overwritten by the value fetched from Orchestrator."
$ws.Range("A20").Value2 = "####  Legend of Key Value pair colours####"
$ws.Range("A19").Value2 = ""

# ---------------------------------------------------------------------------
# 3. Apply styles / row heights.
# ---------------------------------------------------------------------------

# Row 1: unchanged big title style (bold 14, centered, wrap) - leave as is.

# Row 2: blank, takes the same "big title" style as row 1 (general
# alignment, bold 14, wrap text), taller row.
$r = $ws.Range("A2")
$r.Font.Name = "Calibri"
$r.Font.Size = 14
$r.Font.Bold = $true
$r.Font.Color = 0
$r.WrapText = $true
$ws.Rows.Item(2).RowHeight = 18.75

# Row 3: plain wrap-text paragraph style (general alignment, not centered).
$r = $ws.Range("A3")
$r.Font.Name = "Calibri"
$r.Font.Size = 11
$r.Font.Bold = $false
$r.Font.Color = 0
$r.WrapText = $true
$ws.Rows.Item(3).RowHeight = 60

function Set-HeadingStyle($range) {
    $range.Font.Name = "Calibri"
    $range.Font.Size = 12
    $range.Font.Bold = $true
    $range.Font.Color = 0
    $range.HorizontalAlignment = -4108   # xlCenter
    $range.WrapText = $true
}

function Set-DescriptionStyle($range) {
    $range.Font.Name = "Calibri"
    $range.Font.Size = 11
    $range.Font.Bold = $false
    $range.Font.Color = 0
    $range.HorizontalAlignment = -4131   # xlLeft
    $range.WrapText = $true
}

$headingRows = @(4, 6, 8, 10, 15, 17)
foreach ($rowNum in $headingRows) {
    Set-HeadingStyle ($ws.Range("A$rowNum"))
    $ws.Rows.Item($rowNum).RowHeight = 15.75
}

$descRows = @(5, 7, 9, 11, 12, 13, 14, 16, 18, 19)
foreach ($rowNum in $descRows) {
    Set-DescriptionStyle ($ws.Range("A$rowNum"))
}

$ws.Rows.Item(7).RowHeight = 45
$ws.Rows.Item(11).RowHeight = 30
$ws.Rows.Item(12).RowHeight = 45
$ws.Rows.Item(13).RowHeight = 45
$ws.Rows.Item(16).RowHeight = 45
$ws.Rows.Item(18).RowHeight = 30

# Row 20: old "Legend" row - text changed, style (bold 11, wrap) unchanged.

# Row 21: old row 4 ("You may want to mark keys...") keeps wrap-text but now
# on the family-2 flavour of the plain font (distinct xf from row 1-19's).
$r = $ws.Range("A21")
$r.Font.Name = "Calibri"
$r.Font.Size = 11
$r.Font.Bold = $false
$r.Font.Color = 0
$r.WrapText = $true

# ---------------------------------------------------------------------------
# 4. Column width & selection cosmetics for the Introduction sheet.
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 117.15
$ws.Range("A2").Select()

# ---------------------------------------------------------------------------
# 5. Cosmetic selection changes on other sheets (cursor position only).
# ---------------------------------------------------------------------------
$wsSettings = $wb.Worksheets.Item("Settings")
$wsSettings.Range("A8").Select()

$wsTasks = $wb.Worksheets.Item("Tasks")
$wsTasks.Range("B4").Select()

$wsConstants = $wb.Worksheets.Item("Constants")
$wsConstants.Range("B3").Select()

$ws.Activate()
